$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 74: 'Adhesive of Antipathy' / 'Wing Glue'
$ws.Range("I74").Value = 4000
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 4000
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = -3064
$ws.Range("N74").ClearContents()

# Row 77: "It's Gonna Grow Back (L)" / 'Wing Glue'
$ws.Range("I77").Value = 4000
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 20000
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = -15320
$ws.Range("N77").ClearContents()

# Row 98: 'The Dotted Line' / 'Enchanted Durium Ink'
$ws.Range("H98").Value = 10850.909
$ws.Range("I98").Value = 6936
$ws.Range("J98").Value = 50000
$ws.Range("K98").Value = 6936
$ws.Range("L98").Value = 50000
$ws.Range("M98").Value = -5438
$ws.Range("N98").Value = -52996

# Row 122: 'Wishful Inking' / 'Enchanted High Durium Ink'
$ws.Range("H122").Value = 10850.909
$ws.Range("I122").Value = 6936
$ws.Range("J122").Value = 50000
$ws.Range("K122").Value = 20808
$ws.Range("L122").Value = 150000
$ws.Range("M122").Value = -18358
$ws.Range("N122").Value = -154900

# Row 132: 'Fast-forwarding Flora' / 'Growth Formula Lambda'
$ws.Range("H132").Value = 12827782
$ws.Range("I132").Value = 16675127
$ws.Range("J132").Value = 3300.1667
$ws.Range("K132").Value = 50025381
$ws.Range("L132").Value = 9900.500100000001
$ws.Range("M132").Value = -50022851
$ws.Range("N132").Value = -14960.5001

# Row 133: 'Big Brush, Big Dreams' / 'Ginseng Angle Brush'
$ws.Range("H133").Value = 29343.8
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 29343.8
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 29343.8
$ws.Range("N133").Value = -39463.8

# Row 137: 'Cutting Edge of Culinary Quality' / 'Magnesia Whetstone'
$ws.Range("H137").Value = 1168.6052
$ws.Range("I137").Value = 799.4400000000001
$ws.Range("J137").Value = 1878.5385
$ws.Range("K137").Value = 2398.32
$ws.Range("L137").Value = 5635.6155
$ws.Range("M137").Value = 151.6799999999998
$ws.Range("N137").Value = -10735.6155

# Row 138: 'All-night Crafting' / "Cunning Craftsman's Tisane"
$ws.Range("H138").Value = 670420.9
$ws.Range("I138").Value = 1151.9656
$ws.Range("J138").Value = 1209554.1
$ws.Range("K138").Value = 3455.8968
$ws.Range("L138").Value = 3628662.3
$ws.Range("M138").Value = 1684.1032
$ws.Range("N138").Value = -3638942.3

# Row 141: 'Remedy for Reason' / 'Grade 1 Gemdraught of Mind'
$ws.Range("H141").Value = 1072.7778
$ws.Range("I141").Value = 777.8570999999999
$ws.Range("J141").Value = 2105
$ws.Range("K141").Value = 2333.5713
$ws.Range("L141").Value = 6315
$ws.Range("M141").Value = 2846.4287
$ws.Range("N141").Value = -16675

$ws = $wb.Worksheets.Item("ARM")
# Row 32: 'Ingot We Trust' / 'Steel Ingot'
$ws.Range("H32").Value = 3344.2236
$ws.Range("I32").Value = 3012.4478
$ws.Range("J32").Value = 5814.1113
$ws.Range("K32").Value = 3012.4478
$ws.Range("L32").Value = 5814.1113
$ws.Range("M32").Value = -2725.4478
$ws.Range("N32").Value = -6388.1113

# Row 74: 'As the Bolt Flies' / 'Titanium Nugget'
$ws.Range("H74").Value = 1193.4
$ws.Range("I74").Value = 848
$ws.Range("J74").Value = 2575
$ws.Range("K74").Value = 848
$ws.Range("L74").Value = 2575
$ws.Range("M74").Value = 26
$ws.Range("N74").Value = -4323

# Row 77: 'Heavy Metal Banned (L)' / 'Titanium Nugget'
$ws.Range("H77").Value = 1193.4
$ws.Range("I77").Value = 848
$ws.Range("J77").Value = 2575
$ws.Range("K77").Value = 4240
$ws.Range("L77").Value = 12875
$ws.Range("M77").Value = 128
$ws.Range("N77").Value = -21611

# Row 107: 'Shielding the Realm' / 'Deepgold Kite Shield'
$ws.Range("H107").Value = 0
$ws.Range("I107").Value = 0
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 0
$ws.Range("L107").Value = 0
$ws.Range("N107").ClearContents()

# Row 132: "Don't Bore Me, Ore Me" / 'Mountain Chromite Ingot'
$ws.Range("H132").Value = 2394.4375
$ws.Range("I132").Value = 2126.0833
$ws.Range("J132").Value = 3199.5
$ws.Range("K132").Value = 6378.249899999999
$ws.Range("L132").Value = 9598.5
$ws.Range("M132").Value = -3848.249899999999
$ws.Range("N132").Value = -14658.5

$ws = $wb.Worksheets.Item("BSM")
# Row 124: 'History of the Hrothgar' / 'High Durium Bayonet'
$ws.Range("H124").Value = 40624
$ws.Range("I124").Value = 0
$ws.Range("J124").Value = 40624
$ws.Range("K124").Value = 0
$ws.Range("L124").Value = 40624
$ws.Range("N124").Value = -50444

$ws = $wb.Worksheets.Item("CRP")
# Row 31: 'Wall Not Found' / 'Walnut Lumber'
$ws.Range("H31").Value = 1224.4138
$ws.Range("I31").Value = 1202.6923
$ws.Range("J31").Value = 1412.6666
$ws.Range("K31").Value = 1202.6923
$ws.Range("L31").Value = 1412.6666
$ws.Range("M31").Value = -907.6922999999999
$ws.Range("N31").Value = -2002.6666

# Row 34: 'Armoires of the Rich and Famous' / 'Walnut Lumber'
$ws.Range("H34").Value = 1224.4138
$ws.Range("I34").Value = 1202.6923
$ws.Range("J34").Value = 1412.6666
$ws.Range("K34").Value = 1202.6923
$ws.Range("L34").Value = 1412.6666
$ws.Range("M34").Value = -1000.6923
$ws.Range("N34").Value = -1816.6666

# Row 99: 'O Pine' / 'Pine Lumber'
$ws.Range("H99").Value = 2001
$ws.Range("I99").Value = 2002.4
$ws.Range("J99").Value = 1997.5
$ws.Range("K99").Value = 2002.4
$ws.Range("L99").Value = 1997.5
$ws.Range("M99").Value = -504.4000000000001
$ws.Range("N99").Value = -4993.5

# Row 112: 'Understaffed' / 'Applewood Cane'
$ws.Range("H112").Value = 36513
$ws.Range("I112").Value = 0
$ws.Range("J112").Value = 36513
$ws.Range("K112").Value = 0
$ws.Range("L112").Value = 36513
$ws.Range("N112").Value = -39467

# Row 126: 'A Better Conductor' / 'Red Pine Lumber'
$ws.Range("H126").Value = 2001
$ws.Range("I126").Value = 2002.4
$ws.Range("J126").Value = 1997.5
$ws.Range("K126").Value = 6007.200000000001
$ws.Range("L126").Value = 5992.5
$ws.Range("M126").Value = -3537.200000000001
$ws.Range("N126").Value = -10932.5

# Row 132: 'Hull Lotta Damage' / 'Ginseng Lumber'
$ws.Range("H132").Value = 6176.3076
$ws.Range("I132").Value = 7049.263
$ws.Range("J132").Value = 3806.8572
$ws.Range("K132").Value = 21147.789
$ws.Range("L132").Value = 11420.5716
$ws.Range("M132").Value = -18617.789
$ws.Range("N132").Value = -16480.5716

$ws = $wb.Worksheets.Item("CUL")
# Row 107: 'Slippery Service' / 'Frantoio Oil'
$ws.Range("H107").Value = 9552.362999999999
$ws.Range("I107").Value = 480
$ws.Range("J107").Value = 12954.5
$ws.Range("K107").Value = 1440
$ws.Range("L107").Value = 38863.5
$ws.Range("M107").Value = 480
$ws.Range("N107").Value = -42703.5

# Row 132: 'More Mezcal' / 'Cooking Mezcal'
$ws.Range("H132").Value = 1781.25
$ws.Range("I132").Value = 1160
$ws.Range("J132").Value = 2816.6667
$ws.Range("K132").Value = 10440
$ws.Range("L132").Value = 25350.0003
$ws.Range("M132").Value = -7910
$ws.Range("N132").Value = -30410.0003

$ws = $wb.Worksheets.Item("GSM")
# Row 132: 'On Board for Lar' / 'Lar Ingot'
$ws.Range("H132").Value = 2300.5
$ws.Range("I132").Value = 1874.6522
$ws.Range("J132").Value = 2953.4666
$ws.Range("K132").Value = 5623.9566
$ws.Range("L132").Value = 8860.399800000001
$ws.Range("M132").Value = -3093.9566
$ws.Range("N132").Value = -13920.3998

# Row 140: 'The Right Rod' / "Ra'Kaznar Rod"
$ws.Range("H140").Value = 29416.6
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 29416.6
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 29416.6
$ws.Range("N140").Value = -39776.6

$ws = $wb.Worksheets.Item("LTW")
# Row 7: 'Tan Before the Ban' / 'Leather'
$ws.Range("H7").Value = 2136.3635
$ws.Range("I7").Value = 2025
$ws.Range("J7").Value = 2433.3333
$ws.Range("K7").Value = 2025
$ws.Range("L7").Value = 2433.3333
$ws.Range("M7").Value = -1913
$ws.Range("N7").Value = -2657.3333

# Row 44: 'The Righteous Tools for the Job' / 'Boarskin Ringbands'
$ws.Range("H44").Value = 10500
$ws.Range("I44").Value = 0
$ws.Range("J44").Value = 10500
$ws.Range("K44").Value = 0
$ws.Range("L44").Value = 10500
$ws.Range("N44").Value = -11412

# Row 46: 'Supply Side Logic' / 'Boar Leather'
$ws.Range("H46").Value = 6350
$ws.Range("I46").Value = 1400
$ws.Range("J46").Value = 8000
$ws.Range("K46").Value = 1400
$ws.Range("L46").Value = 8000
$ws.Range("M46").Value = -1212
$ws.Range("N46").Value = -8376

# Row 61: 'Spelling Me Softly' / 'Raptor Leather'
$ws.Range("H61").Value = 2061.8
$ws.Range("I61").Value = 1781.6
$ws.Range("J61").Value = 2342
$ws.Range("K61").Value = 1781.6
$ws.Range("L61").Value = 2342
$ws.Range("M61").Value = -1579.6
$ws.Range("N61").Value = -2746

# Row 101: 'A Stitch in Time' / 'Marid Leather Gloves of Healing'
$ws.Range("H101").Value = 13500
$ws.Range("I101").Value = 0
$ws.Range("J101").Value = 13500
$ws.Range("K101").Value = 0
$ws.Range("L101").Value = 13500
$ws.Range("N101").Value = -19990

# Row 104: 'Brace Yourselves' / 'Gazelleskin Bracers of Fending'
$ws.Range("H104").Value = 6047.125
$ws.Range("I104").Value = 0
$ws.Range("J104").Value = 6047.125
$ws.Range("K104").Value = 0
$ws.Range("L104").Value = 6047.125
$ws.Range("N104").Value = -13035.125

# Row 113: 'Peace in Rest' / 'Atrociraptor Leather'
$ws.Range("H113").Value = 2061.8
$ws.Range("I113").Value = 1781.6
$ws.Range("J113").Value = 2342
$ws.Range("K113").Value = 1781.6
$ws.Range("L113").Value = 2342
$ws.Range("M113").Value = 388.4000000000001
$ws.Range("N113").Value = -6682

# Row 126: 'Battered Books' / 'Saiga Leather'
$ws.Range("H126").Value = 2136.3635
$ws.Range("I126").Value = 2025
$ws.Range("J126").Value = 2433.3333
$ws.Range("K126").Value = 6075
$ws.Range("L126").Value = 7299.999899999999
$ws.Range("M126").Value = -3605
$ws.Range("N126").Value = -12239.9999

# Row 132: 'Tenets of Tanning' / 'Silver Lobo Leather'
$ws.Range("H132").Value = 21674.3
$ws.Range("I132").Value = 1296.9375
$ws.Range("J132").Value = 57900.723
$ws.Range("K132").Value = 3890.8125
$ws.Range("L132").Value = 173702.169
$ws.Range("M132").Value = -1360.8125
$ws.Range("N132").Value = -178762.169

# Row 133: 'The Perfect Accessory' / 'Loboskin Amulet of Fending'
$ws.Range("H133").Value = 35383
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 35383
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 35383
$ws.Range("N133").Value = -40443

$ws = $wb.Worksheets.Item("WVR")
# Row 112: 'Hair Do No Harm' / 'Iridescent Hat of Healing'
$ws.Range("H112").Value = 0
$ws.Range("I112").Value = 0
$ws.Range("J112").Value = 0
$ws.Range("K112").Value = 0
$ws.Range("L112").Value = 0
$ws.Range("N112").ClearContents()

# Row 132: 'Comfy Cabins' / 'Snow Cotton Cloth'
$ws.Range("H132").Value = 2742.5366
$ws.Range("I132").Value = 2801.7334
$ws.Range("J132").Value = 2581.0908
$ws.Range("K132").Value = 8405.200199999999
$ws.Range("L132").Value = 7743.2724
$ws.Range("M132").Value = -5875.200199999999
$ws.Range("N132").Value = -12803.2724
